$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 61.625
$ws.Range("I9").Value = 88
$ws.Range("J9").Value = 17.666666
$ws.Range("K9").Value = 88
$ws.Range("L9").Value = 17.666666
$ws.Range("M9").Value = 81
$ws.Range("N9").Value = -355.666666

$ws.Range("H49").Value = 784
$ws.Range("I49").Value = 784
$ws.Range("K49").Value = 2352
$ws.Range("M49").Value = -2216

$ws.Range("H106").Value = 10029.286
$ws.Range("I106").Value = 7441.2
$ws.Range("K106").Value = 7441.2
$ws.Range("M106").Value = -6810.2

$ws.Range("H107").Value = 1521.4546
$ws.Range("I107").Value = 1204
$ws.Range("J107").Value = 2950
$ws.Range("K107").Value = 1204
$ws.Range("L107").Value = 2950
$ws.Range("M107").Value = 716
$ws.Range("N107").Value = -6790

$ws.Range("H111").Value = 3049.5
$ws.Range("I111").Value = 2100
$ws.Range("J111").Value = 3999
$ws.Range("K111").Value = 6300
$ws.Range("L111").Value = 11997
$ws.Range("M111").Value = -3233
$ws.Range("N111").Value = -18131

$ws.Range("H113").Value = 4214.1
$ws.Range("I113").Value = 3269.4
$ws.Range("J113").Value = 5158.8
$ws.Range("K113").Value = 3269.4
$ws.Range("L113").Value = 5158.8
$ws.Range("M113").Value = -15.40000000000009
$ws.Range("N113").Value = -11666.8

$ws.Range("H127").Value = 3079.6
$ws.Range("I127").Value = 3032.6667
$ws.Range("J127").Value = 3150
$ws.Range("K127").Value = 9098.000100000001
$ws.Range("L127").Value = 9450
$ws.Range("M127").Value = -4138.000100000001
$ws.Range("N127").Value = -19370

$ws.Range("H141").Value = 5373.846
$ws.Range("I141").Value = 6096.4546
$ws.Range("K141").Value = 18289.3638
$ws.Range("M141").Value = -13109.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24219.82
$ws.Range("I32").Value = 26841.365
$ws.Range("K32").Value = 26841.365
$ws.Range("M32").Value = -26554.365

$ws.Range("H41").Value = 2171.6
$ws.Range("I41").Value = 1286
$ws.Range("J41").Value = 3500
$ws.Range("K41").Value = 1286
$ws.Range("L41").Value = 3500
$ws.Range("M41").Value = -872
$ws.Range("N41").Value = -4328

$ws.Range("H61").Value = 3915.3333
$ws.Range("I61").Value = 3915.3333
$ws.Range("K61").Value = 3915.3333
$ws.Range("M61").Value = -3703.3333

$ws.Range("H97").Value = 1946.8235
$ws.Range("I97").Value = 1303.08
$ws.Range("J97").Value = 3735
$ws.Range("K97").Value = 1303.08
$ws.Range("L97").Value = 3735
$ws.Range("M97").Value = -807.0799999999999
$ws.Range("N97").Value = -4727

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 2873
$ws.Range("I122").Value = 1996
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 5988
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -3538
$ws.Range("N122").Value = -16150

$ws.Range("H132").Value = 20205.264
$ws.Range("I132").Value = 21352.52
$ws.Range("K132").Value = 64057.56
$ws.Range("M132").Value = -61527.56

$ws.Range("H136").Value = 3915.3333
$ws.Range("I136").Value = 3915.3333
$ws.Range("K136").Value = 11745.9999
$ws.Range("M136").Value = -9195.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 324.1111
$ws.Range("J80").Value = 402.66666
$ws.Range("L80").Value = 402.66666
$ws.Range("N80").Value = -2398.66666

$ws.Range("H83").Value = 324.1111
$ws.Range("J83").Value = 402.66666
$ws.Range("L83").Value = 2013.3333
$ws.Range("N83").Value = -11997.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 5000
$ws.Range("L32").Value = 5000
$ws.Range("N32").Value = -5632

$ws.Range("H58").Value = 127669.75
$ws.Range("I58").Value = 145336.86
$ws.Range("K58").Value = 145336.86
$ws.Range("M58").Value = -145133.86

$ws.Range("H107").Value = 426.0909
$ws.Range("I107").Value = 426.0909
$ws.Range("K107").Value = 426.0909
$ws.Range("M107").Value = 1493.9091

$ws.Range("H136").Value = 127669.75
$ws.Range("I136").Value = 145336.86
$ws.Range("K136").Value = 436010.58
$ws.Range("M136").Value = -433460.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 409
$ws.Range("J46").Value = 450
$ws.Range("L46").Value = 1350
$ws.Range("N46").Value = -1532

$ws.Range("H114").Value = 873.8889
$ws.Range("I114").Value = 504
$ws.Range("K114").Value = 1512
$ws.Range("M114").Value = 1742

$ws.Range("H117").Value = 7622.25
$ws.Range("I117").Value = 6496.6665
$ws.Range("J117").Value = 10999
$ws.Range("K117").Value = 19489.9995
$ws.Range("L117").Value = 32997
$ws.Range("M117").Value = -16047.9995
$ws.Range("N117").Value = -39881

$ws.Range("H121").Value = 1863.6
$ws.Range("I121").Value = 5100
$ws.Range("J121").Value = 1054.5
$ws.Range("K121").Value = 15300
$ws.Range("L121").Value = 3163.5
$ws.Range("M121").Value = -13990
$ws.Range("N121").Value = -5783.5

$ws.Range("H131").Value = 2132822.5
$ws.Range("I131").Value = 2331.4285
$ws.Range("K131").Value = 6994.2855
$ws.Range("M131").Value = -1954.2855

$ws.Range("H132").Value = 1632.3334
$ws.Range("I132").Value = 1632.3334
$ws.Range("K132").Value = 14691.0006
$ws.Range("M132").Value = -12161.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3276.625
$ws.Range("I102").Value = 1702.3334
$ws.Range("K102").Value = 1702.3334
$ws.Range("M102").Value = -80.33339999999998

$ws.Range("H132").Value = 32524.094
$ws.Range("I132").Value = 38320.074
$ws.Range("K132").Value = 114960.222
$ws.Range("M132").Value = -112430.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2475.5
$ws.Range("I22").Value = 832.5
$ws.Range("J22").Value = 4529.25
$ws.Range("K22").Value = 832.5
$ws.Range("L22").Value = 4529.25
$ws.Range("M22").Value = -537.5
$ws.Range("N22").Value = -5119.25

$ws.Range("H27").Value = 2475.5
$ws.Range("I27").Value = 832.5
$ws.Range("J27").Value = 4529.25
$ws.Range("K27").Value = 832.5
$ws.Range("L27").Value = 4529.25
$ws.Range("M27").Value = -725.5
$ws.Range("N27").Value = -4743.25

$ws.Range("H132").Value = 112643.27
$ws.Range("I132").Value = 136409.22
$ws.Range("K132").Value = 409227.66
$ws.Range("M132").Value = -406697.66

$ws.Range("H136").Value = 3332.625
$ws.Range("I136").Value = 2717.9443
$ws.Range("J136").Value = 5176.6665
$ws.Range("K136").Value = 8153.8329
$ws.Range("L136").Value = 15529.9995
$ws.Range("M136").Value = -5603.8329
$ws.Range("N136").Value = -20629.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1295.48
$ws.Range("I122").Value = 1281.2727
$ws.Range("J122").Value = 1399.6666
$ws.Range("K122").Value = 3843.8181
$ws.Range("L122").Value = 4198.9998
$ws.Range("M122").Value = -1393.8181
$ws.Range("N122").Value = -9098.9998

$ws.Range("H132").Value = 27727.041
$ws.Range("I132").Value = 33468.285
$ws.Range("K132").Value = 100404.855
$ws.Range("M132").Value = -97874.85500000001

$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -100120
